$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (2,3,4,6,7) got shuffled around (a cyclic re-ordering of the
# weekly records), while row 5 stayed put. Columns A,B,C,E,F,G,H,I,J,Q,T
# are identical across all these rows, so only D,K,L,M,N,O,P,R,S need to
# change per row.

# Row 2 -> becomes old row 3 (Hachiya / Segunda)
$ws.Range("D2").Value = 44301
$ws.Range("K2").Value = "Hachiya"
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 250
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1139

# Row 3 -> becomes old row 7 (Mankaki / Primera, 44699)
$ws.Range("D3").Value = 44699
$ws.Range("K3").Value = "Mankaki"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 29000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 29500
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1639

# Row 4 -> becomes old row 2 (Mankaki / Primera, 44313)
$ws.Range("D4").Value = 44313
$ws.Range("K4").Value = "Mankaki"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1194

# Row 5 stays the same (no change per diff)

# Row 6 -> becomes old row 4 (Mankaki / Segunda, 44305)
$ws.Range("D6").Value = 44305
$ws.Range("K6").Value = "Mankaki"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1361

# Row 7 -> becomes old row 6 (Mankaki / Segunda, 44355, Región Metropolitana)
$ws.Range("D7").Value = 44355
$ws.Range("K7").Value = "Mankaki"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 270
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 20500
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 1139
